# Selenium "Upload file" changed logic -> rename worksheets and refresh the
# "summary" label used by the Report test-data rows, and leave the Report
# sheet/tab as the active selection.

$wb = $excel.ActiveWorkbook

$wsReport = $wb.Worksheets.Item(1)
$wsAudit  = $wb.Worksheets.Item(2)

# Rename sheets: Sheet1 -> Report, Sheet2 -> Audit
$wsReport.Name = "Report"
$wsAudit.Name  = "Audit"

# The "Report" sheet's B column used to read "MATTER SUMMARY" for every row;
# it now reads "CASE SUMMARY".
$wsReport.Range("B2:B6").Value = "CASE SUMMARY"

# Make the Report sheet the active tab, with B2:B6 selected.
$wsReport.Activate()
$wsReport.Range("B2:B6").Select()
